$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# Data for columns I (Minute4), J (Second4), K (Rep4) for rows 2-19
$data = @{
    2  = 9
    3  = 14
    4  = 9
    5  = 12
    6  = 7
    7  = 13
    8  = 13
    9  = 17
    10 = 15
    11 = 13
    12 = 9
    13 = 12
    14 = 20
    15 = 5
    16 = 12
    17 = 18
    18 = 19
    19 = 6
}

foreach ($row in 2..19) {
    $ws.Cells.Item($row, 9).Value = 10
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = $data[$row]
}

$ws.Range("K20").Select()
